$d = $word.ActiveDocument

# Locate the existing "מנחה ערב שבת" paragraph; the new "מנחה גדולה ערב שבת"
# paragraph must be inserted immediately before it.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("מנחה ערב שבת")) {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not locate the 'מנחה ערב שבת' paragraph"
}

$target = $d.Paragraphs.Item($targetIndex)
$target.Range.InsertParagraphBefore()
$newp = $d.Paragraphs.Item($targetIndex)
$xmlFrag = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="List"/>
    <w:spacing w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma" w:hint="cs"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="864904" w:themeColor="accent3" w:themeShade="80"/>
      <w:sz w:val="36"/>
      <w:szCs w:val="36"/>
      <w:rtl/>
      <w:lang w:val="en-GB"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma" w:hint="cs"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="864904" w:themeColor="accent3" w:themeShade="80"/>
      <w:sz w:val="36"/>
      <w:szCs w:val="36"/>
      <w:rtl/>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t xml:space="preserve">מנחה גדולה ערב שבת </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="864904" w:themeColor="accent3" w:themeShade="80"/>
      <w:sz w:val="36"/>
      <w:szCs w:val="36"/>
      <w:rtl/>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t>&#8211;</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma" w:hint="cs"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="864904" w:themeColor="accent3" w:themeShade="80"/>
      <w:sz w:val="36"/>
      <w:szCs w:val="36"/>
      <w:rtl/>
      <w:lang w:val="en-GB"/>
    </w:rPr>
    <w:t xml:space="preserve"> 14:30</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$newp.Range.InsertXML($xmlFrag) | Out-Null
